$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @('Obafemi Martins', 'Nigeria', 'West', 'Forward', 14.3, '03_04', 'Inter', 0.49, 0.42, 0.91, 0.49, 0.91),
    @('Yakubu Aiyegbeni', 'Nigeria', 'West', 'Forward', 12.2, '09_10', 'Everton', 0.41, 0.33, 0.74, 0.41, 0.74),
    @('Kelechi Iheanacho', 'Nigeria', 'West', 'Forward / Midfielder', 6.0, '16_17', 'Manchester City', 0.67, 0.5, 1.17, 0.67, 1.17),
    @('Ademola Lookman', 'Nigeria', 'West', 'Forward / Midfielder', 6.4, '17_18', 'RB Leipzig', 0.78, 0.47, 1.25, 0.78, 1.25),
    @('Sébastien Haller', 'Côte d’Ivoire', 'West', 'Forward', 13.6, '22_23', 'Dortmund', 0.66, 0.22, 0.88, 0.66, 0.88),
    @('Inaki Williams', 'Ghana', 'West', 'Forward / Midfielder', 18.0, '15_16', 'Athletic Club', 0.44, 0.17, 0.61, 0.44, 0.61),
    @('Paul Onuachu', 'Nigeria', 'West', 'Forward', 11.7, '24_25', 'Southampton', 0.34, 0.09, 0.43, 0.34, 0.43),
    @('Andy Delort', 'Algeria', 'North', 'Forward', 27.6, '20_21', 'Montpellier', 0.54, 0.33, 0.87, 0.54, 0.87),
    @('Youssef En-Nesyri', 'Morocco', 'North', 'Forward', 25.7, '20_21', 'Sevilla', 0.7, 0.0, 0.7, 0.7, 0.7),
    @('Boulaye Dia', 'Senegal', 'West', 'Forward / Midfielder', 28.0, '22_23', 'Salernitana', 0.57, 0.21, 0.79, 0.54, 0.75),
    @('Simon Adingra', 'Côte d’Ivoire', 'West', 'Defender / Forward / Midfielder', 12.2, '24_25', 'Brighton', 0.16, 0.16, 0.33, 0.16, 0.33),
    @('Alex Iwobi', 'Nigeria', 'West', 'Forward / Midfielder', 7.0, '15_16', 'Arsenal', 0.29, 0.29, 0.57, 0.29, 0.57),
    @('Wilfrid Zaha', 'Côte d’Ivoire', 'West', 'Forward / Midfielder', 33.5, '16_17', 'Crystal Palace', 0.21, 0.24, 0.45, 0.21, 0.45),
    @('Habib Diallo', 'Senegal', 'West', 'Forward', 17.1, '21_22', 'Strasbourg', 0.64, 0.06, 0.7, 0.53, 0.58),
    @('Sofiane Boufal', 'Morocco', 'North', 'Forward / Midfielder', 11.3, '14_15', 'Lille', 0.27, 0.53, 0.8, 0.09, 0.62),
    @('Patson Daka', 'Zambia', 'South', 'Forward', 12.9, '21_22', 'Leicester City', 0.39, 0.23, 0.62, 0.39, 0.62),
    @('Umar Sadiq', 'Nigeria', 'West', 'Forward', 9.2, '24_25', 'Valencia', 0.55, 0.0, 0.55, 0.55, 0.55),
    @('Taiwo Awoniyi', 'Nigeria', 'West', 'Forward', 11.6, '23_24', 'Nottingham Forest', 0.52, 0.26, 0.78, 0.52, 0.78),
    @('Musa Barrow', 'Gambia', 'West', 'Forward / Midfielder', 4.9, '17_18', 'Atalanta', 0.61, 0.2, 0.81, 0.61, 0.81),
    @('Odion Ighalo', 'Nigeria', 'West', 'Forward', 35.0, '15_16', 'Watford', 0.43, 0.09, 0.51, 0.43, 0.51),
    @('Ihlas Bebou', 'Togo', 'West', 'Forward / Midfielder', 18.2, '23_24', 'Hoffenheim', 0.39, 0.28, 0.66, 0.39, 0.66),
    @('Amad Diallo', 'Côte d’Ivoire', 'West', 'Forward / Midfielder', 21.1, '24_25', 'Manchester United', 0.38, 0.28, 0.66, 0.38, 0.66)
)

$startRow = 43
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $rowVals = $newRows[$i]
    for ($j = 0; $j -lt $rowVals.Count; $j++) {
        $ws.Cells.Item($r, $j + 1).Value = $rowVals[$j]
    }
}

$ws.Columns.Item(4).ColumnWidth = 28.08984375

$ws.Range("H64:L64").Select()
